$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "End Time" column (E) values from 100 to 80 for data rows 2-55
$ws.Range("E2:E55").Value = 80

# Update the active selection to match the saved workbook state
$ws.Range("E4").Select()
